# Update "想去人数" (want-to-go count) values in the 展览 and 全部类型 sheets.
# Mapping of row -> new value is the same on both sheets (they list the same
# exhibition events), only the row numbers differ because 全部类型 interleaves
# additional rows from the 演出 (show) sheet.

$wb = $excel.ActiveWorkbook

$updates_zhanlan = @{
    2  = 4685
    3  = 2544
    7  = 139
    8  = 207
    10 = 1778
    11 = 324
    12 = 4022
    13 = 43
    14 = 276
}

$updates_quanbu = @{
    2  = 4685
    3  = 2544
    9  = 139
    10 = 207
    14 = 1778
    15 = 324
    16 = 4022
    17 = 43
    18 = 276
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates_zhanlan.Keys) {
    $ws1.Range("F$row").Value = $updates_zhanlan[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates_quanbu.Keys) {
    $ws4.Range("F$row").Value = $updates_quanbu[$row]
}
